$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1332.3
$ws.Range("I28").Value = 1360.8
$ws.Range("K28").Value = 1360.8
$ws.Range("M28").Value = -875.8

$ws.Range("H100").Value = 5505.5386
$ws.Range("I100").Value = 3872.125
$ws.Range("J100").Value = 8119
$ws.Range("K100").Value = 3872.125
$ws.Range("L100").Value = 8119
$ws.Range("M100").Value = -3331.125
$ws.Range("N100").Value = -9201

$ws.Range("H116").Value = 6416.3335
$ws.Range("I116").Value = 5999.6665
$ws.Range("J116").Value = 6833
$ws.Range("K116").Value = 5999.6665
$ws.Range("L116").Value = 6833
$ws.Range("M116").Value = -2557.6665
$ws.Range("N116").Value = -13717

$ws.Range("H125").Value = 1513.3334
$ws.Range("I125").Value = 1820
$ws.Range("J125").Value = 900
$ws.Range("K125").Value = 16380
$ws.Range("L125").Value = 8100
$ws.Range("M125").Value = -13920
$ws.Range("N125").Value = -13020

$ws.Range("H132").Value = 15893.45
$ws.Range("I132").Value = 2902.4517
$ws.Range("K132").Value = 8707.355100000001
$ws.Range("M132").Value = -6177.355100000001

$ws.Range("H137").Value = 3021.9048
$ws.Range("I137").Value = 2520.9443
$ws.Range("J137").Value = 6027.6665
$ws.Range("K137").Value = 7562.8329
$ws.Range("L137").Value = 18082.9995
$ws.Range("M137").Value = -5012.8329
$ws.Range("N137").Value = -23182.9995

$ws.Range("H138").Value = 2695.1714
$ws.Range("I138").Value = 1336.0571
$ws.Range("K138").Value = 4008.1713
$ws.Range("M138").Value = 1131.8287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20001774
$ws.Range("I32").Value = 21278394
$ws.Range("K32").Value = 21278394
$ws.Range("M32").Value = -21278107

$ws.Range("H61").Value = 794.3714
$ws.Range("I61").Value = 781.9394
$ws.Range("J61").Value = 999.5
$ws.Range("K61").Value = 781.9394
$ws.Range("L61").Value = 999.5
$ws.Range("M61").Value = -569.9394
$ws.Range("N61").Value = -1423.5

$ws.Range("H64").Value = 42500
$ws.Range("J64").Value = 42500
$ws.Range("L64").Value = 42500
$ws.Range("N64").Value = -42996

$ws.Range("H67").Value = 42500
$ws.Range("J67").Value = 42500
$ws.Range("L67").Value = 42500
$ws.Range("N67").Value = -44216

$ws.Range("H74").Value = 1723.3226
$ws.Range("I74").Value = 1814.3846
$ws.Range("K74").Value = 1814.3846
$ws.Range("M74").Value = -940.3846000000001

$ws.Range("H77").Value = 1723.3226
$ws.Range("I77").Value = 1814.3846
$ws.Range("K77").Value = 9071.923000000001
$ws.Range("M77").Value = -4703.923000000001

$ws.Range("H102").Value = 6772.1113
$ws.Range("I102").Value = 6868.625
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 6868.625
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -5246.625
$ws.Range("N102").Value = -9244

$ws.Range("H122").Value = 2262.6
$ws.Range("J122").Value = 3650
$ws.Range("L122").Value = 10950
$ws.Range("N122").Value = -15850

$ws.Range("H132").Value = 2507
$ws.Range("I132").Value = 2507
$ws.Range("K132").Value = 7521
$ws.Range("M132").Value = -4991

$ws.Range("H136").Value = 794.3714
$ws.Range("I136").Value = 781.9394
$ws.Range("J136").Value = 999.5
$ws.Range("K136").Value = 2345.8182
$ws.Range("L136").Value = 2998.5
$ws.Range("M136").Value = 204.1818000000003
$ws.Range("N136").Value = -8098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null

$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314

$ws.Range("H134").Value = 1361.5143
$ws.Range("I134").Value = 1195.6765
$ws.Range("K134").Value = 3587.0295
$ws.Range("M134").Value = -1052.0295

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 24296.428
$ws.Range("J86").Value = 15750.556
$ws.Range("L86").Value = 15750.556
$ws.Range("N86").Value = -17996.556

$ws.Range("H89").Value = 24296.428
$ws.Range("J89").Value = 15750.556
$ws.Range("L89").Value = 78752.78
$ws.Range("N89").Value = -89984.78

$ws.Range("H122").Value = 792113.9
$ws.Range("I122").Value = 2556770
$ws.Range("J122").Value = 7822.222
$ws.Range("K122").Value = 7670310
$ws.Range("L122").Value = 23466.666
$ws.Range("M122").Value = -7667860
$ws.Range("N122").Value = -28366.666

$ws.Range("H134").Value = 4129.8945
$ws.Range("I134").Value = 3469.6667
$ws.Range("K134").Value = 10409.0001
$ws.Range("M134").Value = -7874.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7665.8335
$ws.Range("I56").Value = 7665.8335
$ws.Range("K56").Value = 7665.8335
$ws.Range("M56").Value = -7135.8335

$ws.Range("H131").Value = 2145.4517
$ws.Range("I131").Value = 1089.1904
$ws.Range("K131").Value = 3267.5712
$ws.Range("M131").Value = 1772.4288

$ws.Range("H136").Value = 2916
$ws.Range("I136").Value = 1665.6
$ws.Range("K136").Value = 4996.799999999999
$ws.Range("M136").Value = 103.2000000000007

$ws.Range("H138").Value = 5016.136
$ws.Range("I138").Value = 3022.4167
$ws.Range("K138").Value = 9067.250100000001
$ws.Range("M138").Value = -3927.250100000001

$ws.Range("H141").Value = 76926264
$ws.Range("I141").Value = 76926264
$ws.Range("K141").Value = 230778792
$ws.Range("M141").Value = -230773612

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2919.1562
$ws.Range("I102").Value = 2037.7368
$ws.Range("J102").Value = 4207.385
$ws.Range("K102").Value = 2037.7368
$ws.Range("L102").Value = 4207.385
$ws.Range("M102").Value = -415.7367999999999
$ws.Range("N102").Value = -7451.385

$ws.Range("H126").Value = 6218.778
$ws.Range("J126").Value = 6216
$ws.Range("L126").Value = 18648
$ws.Range("N126").Value = -23588

$ws.Range("H132").Value = 1887.4348
$ws.Range("I132").Value = 1791.409
$ws.Range("K132").Value = 5374.227000000001
$ws.Range("M132").Value = -2844.227000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2562.3257
$ws.Range("I132").Value = 2255.7354
$ws.Range("K132").Value = 6767.206200000001
$ws.Range("M132").Value = -4237.206200000001

$ws.Range("H136").Value = 2750.1765
$ws.Range("I136").Value = 2631.963
$ws.Range("K136").Value = 7895.889000000001
$ws.Range("M136").Value = -5345.889000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = $null

$ws.Range("H51").Value = 27498.75
$ws.Range("I51").Value = 20000
$ws.Range("J51").Value = 34997.5
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 34997.5
$ws.Range("M51").Value = -19490
$ws.Range("N51").Value = -36017.5

$ws.Range("H122").Value = 2361.4666
$ws.Range("I122").Value = 2125.5386
$ws.Range("J122").Value = 3895
$ws.Range("K122").Value = 6376.6158
$ws.Range("L122").Value = 11685
$ws.Range("M122").Value = -3926.6158
$ws.Range("N122").Value = -16585

$ws.Range("H124").Value = 213499.6
$ws.Range("J124").Value = 213499.6
$ws.Range("L124").Value = 213499.6
$ws.Range("N124").Value = -223319.6

$ws.Range("H132").Value = 1442.5555
$ws.Range("I132").Value = 1457.5294
$ws.Range("K132").Value = 4372.5882
$ws.Range("M132").Value = -1842.5882
